$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
  @(465, 44539, 8, 39, 258.5006959634122),
  @(466, 44540, 5, 36, 238.6160270431498),
  @(467, 44541, 2, 32, 212.1031351494664),
  @(468, 44542, 1, 28, 185.5902432557831),
  @(469, 44543, 3, 28, 185.5902432557831),
  @(470, 44544, 3, 25, 165.7055743355207),
  @(471, 44545, 1, 23, 152.449128388679),
  @(472, 44546, 11, 26, 172.3337973089415),
  @(473, 44547, 17, 38, 251.8724729899914),
  @(474, 44548, 2, 38, 251.8724729899914),
  @(475, 44550, 9, 46, 304.898256777358),
  @(476, 44551, 10, 53, 351.2958175913038),
  @(477, 44552, 0, 50, 331.4111486710413),
  @(478, 44553, 14, 63, 417.578047325512),
  @(479, 44554, 4, 56, 371.1804865115662),
  @(480, 44555, 22, 61, 404.3216013786704),
  @(481, 44556, 17, 76, 503.7449459799828),
  @(482, 44557, 18, 85, 563.3989527407703),
  @(483, 44558, 3, 78, 517.0013919268245),
  @(484, 44559, 6, 84, 556.7707297673493),
  @(485, 44560, 42, 112, 742.3609730231325),
  @(486, 44561, 30, 138, 914.694770332074),
  @(487, 44562, 45, 161, 1067.143898720753),
  @(488, 44563, 28, 172, 1140.054351428382),
  @(489, 44564, 21, 175, 1159.939020348645),
  @(490, 44565, 6, 178, 1179.823689268907),
  @(491, 44566, 22, 194, 1285.87525684364)
)

foreach ($r in $rows) {
  $rowNum = $r[0]
  $ws.Cells.Item($rowNum, 1).Value = $r[1]
  $ws.Cells.Item($rowNum, 2).Value = $r[2]
  $ws.Cells.Item($rowNum, 3).Value = $r[3]
  $ws.Cells.Item($rowNum, 4).Value = $r[4]
}

# Copy the date-column style (bold, centered, bordered, date-formatted) down to the new rows
$ws.Range("A464").Copy()
$ws.Range("A465:A491").PasteSpecial(-4122)
